# Insert a new "Properties" worksheet between "Metadata" and "Concepts",
# describing the custom FHIR CodeSystem property "inativo" (boolean).

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")

# Create the new sheet right after "Metadata" so the final order is
# Metadata, Properties, Concepts.
$wsProperties = $wb.Worksheets.Add($null, $wsMetadata)
$wsProperties.Name = "Properties"

# Header row
$wsProperties.Range("A1").Value = "Code"
$wsProperties.Range("B1").Value = "Uri"
$wsProperties.Range("C1").Value = "Description"
$wsProperties.Range("D1").Value = "Type"

# Data row describing the "inativo" property
$wsProperties.Range("A2").Value = "inativo"
$wsProperties.Range("D2").Value = "boolean"

# Match the formatting already used on the other sheets: bold header row,
# plain body row (including the empty B2/C2 cells), by copying the
# existing styles from the Metadata sheet.
$wsMetadata.Range("A1:B1").Copy()
$wsProperties.Range("A1:D1").PasteSpecial(-4122)

$wsMetadata.Range("A2:B2").Copy()
$wsProperties.Range("A2:D2").PasteSpecial(-4122)
